$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.822732333489625
$ws.Range("C2").Value = 0.1438761311182475
$ws.Range("D2").Value = 0.08669589847261072
$ws.Range("F2").Value = 2.318823329023772
$ws.Range("G2").Value = 1.667700327157945
$ws.Range("H2").Value = 1.467641665125001
$ws.Range("J2").Value = 0.1843079618869989
$ws.Range("L2").Value = 0.328637806213365
$ws.Range("M2").Value = 0.4129372486134599
$ws.Range("N2").Value = 1.986704101451016
$ws.Range("B3").Value = 1.736387436105417
$ws.Range("C3").Value = 0.1288516908818167
$ws.Range("D3").Value = 0.08662404270775426
$ws.Range("F3").Value = 2.317698775571486
$ws.Range("G3").Value = 1.661306650256165
$ws.Range("H3").Value = 1.470778917695824
$ws.Range("J3").Value = 0.1851841030991634
$ws.Range("L3").Value = 0.3265150969484978
$ws.Range("M3").Value = 0.4000130479750155
$ws.Range("N3").Value = 2.009185763621311
$ws.Range("B4").Value = 1.684193627239608
$ws.Range("C4").Value = 0.119570083756912
$ws.Range("D4").Value = 0.08658673871838651
$ws.Range("F4").Value = 2.318227547712297
$ws.Range("G4").Value = 1.658403170196124
$ws.Range("H4").Value = 1.47340525311057
$ws.Range("J4").Value = 0.1857503598464727
$ws.Range("L4").Value = 0.3253312888878668
$ws.Range("M4").Value = 0.3922635746945886
$ws.Range("N4").Value = 2.023672147662906
$ws.Range("B5").Value = 1.66313167803213
$ws.Range("C5").Value = 0.1157735360718277
$ws.Range("D5").Value = 0.08657325684536943
$ws.Range("F5").Value = 2.318749428164338
$ws.Range("G5").Value = 1.657476619145214
$ws.Range("H5").Value = 1.474651489308741
$ws.Range("J5").Value = 0.1859882451909645
$ws.Range("L5").Value = 0.3248790192520303
$ws.Range("M5").Value = 0.3891525246048602
$ws.Range("N5").Value = 2.029747028582432
$ws.Range("B6").Value = 1.659646907138494
$ws.Range("C6").Value = 0.1151422635634844
$ws.Range("D6").Value = 0.08657112224196339
$ws.Range("F6").Value = 2.318854586480711
$ws.Range("G6").Value = 1.657338254435899
$ws.Range("H6").Value = 1.474869053432428
$ws.Range("J6").Value = 0.1860281769923047
$ws.Range("L6").Value = 0.3248057431209901
$ws.Range("M6").Value = 0.3886387753862195
$ws.Range("N6").Value = 2.030766117819385
$ws.Range("B7").Value = 1.683908737565162
$ws.Range("C7").Value = 0.1195189397509466
$ws.Range("D7").Value = 0.08658654992464854
$ws.Range("F7").Value = 2.318233345623753
$ws.Range("G7").Value = 1.658389635831085
$ws.Range("H7").Value = 1.473421347810358
$ws.Range("J7").Value = 0.1857535391588963
$ws.Range("L7").Value = 0.3253250672695032
$ws.Range("M7").Value = 0.392221427805552
$ws.Range("N7").Value = 2.02375338108489
$ws.Range("B8").Value = 1.792790315327522
$ws.Range("C8").Value = 0.138707434566868
$ws.Range("D8").Value = 0.08666971144677049
$ws.Range("F8").Value = 2.318182405515685
$ws.Range("G8").Value = 1.665283358807201
$ws.Range("H8").Value = 1.468578061009808
$ws.Range("J8").Value = 0.184604189811358
$ws.Range("L8").Value = 0.3278811320816999
$ws.Range("M8").Value = 0.4084424520530305
$ws.Range("N8").Value = 1.994314073533297
$ws.Range("B9").Value = 2.012814785227249
$ws.Range("C9").Value = 0.1758891737414672
$ws.Range("D9").Value = 0.08688664311965333
$ws.Range("F9").Value = 2.327768509018398
$ws.Range("G9").Value = 1.686935070529756
$ws.Range("H9").Value = 1.464638666930483
$ws.Range("J9").Value = 0.1825742230111249
$ws.Range("L9").Value = 0.3338390868497214
$ws.Range("M9").Value = 0.4417238197701252
$ws.Range("N9").Value = 1.942001743273536
$ws.Range("B10").Value = 2.178431262258641
$ws.Range("C10").Value = 0.2029392216658721
$ws.Range("D10").Value = 0.08707859536616347
$ws.Range("F10").Value = 2.340737237314229
$ws.Range("G10").Value = 1.707835737712884
$ws.Range("H10").Value = 1.465140076587403
$ws.Range("J10").Value = 0.1812183541596539
$ws.Range("L10").Value = 0.3387896081291473
$ws.Range("M10").Value = 0.4670705672728701
$ws.Range("N10").Value = 1.906872760855929
$ws.Range("B11").Value = 2.254636069317485
$ws.Range("C11").Value = 0.2151881729815273
$ws.Range("D11").Value = 0.08717294379651541
$ws.Range("F11").Value = 2.34792868323386
$ws.Range("G11").Value = 1.718436091510767
$ws.Range("H11").Value = 1.466107266415236
$ws.Range("J11").Value = 0.1806307645384995
$ws.Range("L11").Value = 0.3411655821145558
$ws.Range("M11").Value = 0.4787954904327307
$ws.Range("N11").Value = 1.891609743998353
$ws.Range("B12").Value = 2.283616885837887
$ws.Range("C12").Value = 0.2198184860869503
$ws.Range("D12").Value = 0.08720967726448237
$ws.Range("F12").Value = 2.350837994936967
$ws.Range("G12").Value = 1.722607834254262
$ws.Range("H12").Value = 1.466579911548109
$ws.Range("J12").Value = 0.180412444155762
$ws.Range("L12").Value = 0.3420830597704878
$ws.Range("M12").Value = 0.4832633015027454
$ws.Range("N12").Value = 1.885933282710448
$ws.Range("B13").Value = 2.277369857195083
$ws.Range("C13").Value = 0.218821625667232
$ws.Range("D13").Value = 0.08720172140004223
$ws.Range("F13").Value = 2.350203142507468
$ws.Range("G13").Value = 1.721702356451914
$ws.Range("H13").Value = 1.466473385390572
$ws.Range("J13").Value = 0.1804592773282012
$ws.Range("L13").Value = 0.3418846764185304
$ws.Range("M13").Value = 0.4822998426041423
$ws.Range("N13").Value = 1.887151212776812
$ws.Range("B14").Value = 2.257017863077749
$ws.Range("C14").Value = 0.2155692739967776
$ws.Range("D14").Value = 0.08717594575189835
$ws.Range("F14").Value = 2.34816430330261
$ws.Range("G14").Value = 1.718776141425707
$ws.Range("H14").Value = 1.466144018332812
$ws.Range("J14").Value = 0.1806127193435376
$ws.Range("L14").Value = 0.3412407082983009
$ws.Range("M14").Value = 0.4791625028395003
$ws.Range("N14").Value = 1.891140666650351
$ws.Range("B15").Value = 2.244567769628645
$ws.Range("C15").Value = 0.2135760596851242
$ws.Range("D15").Value = 0.08716028826070144
$ws.Range("F15").Value = 2.3469396956244
$ws.Range("G15").Value = 1.717004292160453
$ws.Range("H15").Value = 1.465956130112119
$ws.Range("J15").Value = 0.1807072519277662
$ws.Range("L15").Value = 0.3408485682543017
$ws.Range("M15").Value = 0.4772444127208288
$ws.Range("N15").Value = 1.893597781451255
$ws.Range("B16").Value = 2.173468447244431
$ws.Range("C16").Value = 0.2021375910951519
$ws.Range("D16").Value = 0.08707257047909067
$ws.Range("F16").Value = 2.340293280781495
$ws.Range("G16").Value = 1.707165008078903
$ws.Range("H16").Value = 1.465091750315281
$ws.Range("J16").Value = 0.1812573410993368
$ws.Range("L16").Value = 0.3386368203753989
$ws.Range("M16").Value = 0.466308218983869
$ws.Range("N16").Value = 1.907884683064622
$ws.Range("B17").Value = 2.130072385648248
$ws.Range("C17").Value = 0.1951060464290038
$ws.Range("D17").Value = 0.08702055504744877
$ws.Range("F17").Value = 2.336547015567703
$ws.Range("G17").Value = 1.701409114701562
$ws.Range("H17").Value = 1.464750859060217
$ws.Range("J17").Value = 0.1816022736203351
$ws.Range("L17").Value = 0.337311678048053
$ws.Range("M17").Value = 0.4596489495438405
$ws.Range("N17").Value = 1.916833073125185
$ws.Range("B18").Value = 2.105193580750381
$ws.Range("C18").Value = 0.1910563938518237
$ws.Range("D18").Value = 0.08699129897566316
$ws.Range("F18").Value = 2.334513845436817
$ws.Range("G18").Value = 1.698201280714159
$ws.Range("H18").Value = 1.46462434784894
$ws.Range("J18").Value = 0.1818034190241136
$ws.Range("L18").Value = 0.3365611632115986
$ws.Range("M18").Value = 0.4558370351490595
$ws.Range("N18").Value = 1.922047459283976
$ws.Range("B19").Value = 2.096784056983324
$ws.Range("C19").Value = 0.189684342127407
$ws.Range("D19").Value = 0.08698150719506081
$ws.Range("F19").Value = 2.333846322051514
$ws.Range("G19").Value = 1.697132804293858
$ws.Range("H19").Value = 1.46459345709448
$ws.Range("J19").Value = 0.1818719959764215
$ws.Range("L19").Value = 0.3363090588150328
$ws.Range("M19").Value = 0.4545495374326265
$ws.Range("N19").Value = 1.923824550769402
$ws.Range("B20").Value = 2.134683544447682
$ws.Range("C20").Value = 0.1958551143560783
$ws.Range("D20").Value = 0.08702602371894841
$ws.Range("F20").Value = 2.336933226445026
$ws.Range("G20").Value = 1.702011195873013
$ws.Range("H20").Value = 1.464779947476757
$ws.Range("J20").Value = 0.1815652705065007
$ws.Range("L20").Value = 0.3374515342649715
$ws.Range("M20").Value = 0.460355944985082
$ws.Range("N20").Value = 1.915873513659111
$ws.Range("B21").Value = 2.262992385027019
$ws.Range("C21").Value = 0.2165247883991697
$ws.Range("D21").Value = 0.08718348942649179
$ws.Range("F21").Value = 2.348758107633259
$ws.Range("G21").Value = 1.719631359869283
$ws.Range("H21").Value = 1.466237873098294
$ws.Range("J21").Value = 0.1805675361649111
$ws.Range("L21").Value = 0.3414293762833864
$ws.Range("M21").Value = 0.4800832611980113
$ws.Range("N21").Value = 1.889966062181385
$ws.Range("B22").Value = 2.347570395758112
$ws.Range("C22").Value = 0.2299864044246931
$ws.Range("D22").Value = 0.08729226259427847
$ws.Range("F22").Value = 2.35757103223726
$ws.Range("G22").Value = 1.732066119896814
$ws.Range("H22").Value = 1.467810895131436
$ws.Range("J22").Value = 0.1799398589020571
$ws.Range("L22").Value = 0.3441325389529197
$ws.Range("M22").Value = 0.493138394843065
$ws.Range("N22").Value = 1.873636304184729
$ws.Range("B23").Value = 2.302363792525455
$ws.Range("C23").Value = 0.2228060095257263
$ws.Range("D23").Value = 0.08723367366841295
$ws.Range("F23").Value = 2.352768062243129
$ws.Range("G23").Value = 1.725345198446547
$ws.Range("H23").Value = 1.466914559683687
$ws.Range("J23").Value = 0.1802726334681761
$ws.Range("L23").Value = 0.3426803731246935
$ws.Range("M23").Value = 0.4861558297115209
$ws.Range("N23").Value = 1.882296642285689
$ws.Range("B24").Value = 2.132598620929343
$ws.Range("C24").Value = 0.1955164829403486
$ws.Range("D24").Value = 0.0870235493135656
$ws.Range("F24").Value = 2.336758244830889
$ws.Range("G24").Value = 1.701738679444759
$ws.Range("H24").Value = 1.464766580199125
$ws.Range("J24").Value = 0.1815819907695708
$ws.Range("L24").Value = 0.3373882699809627
$ws.Range("M24").Value = 0.4600362606639976
$ws.Range("N24").Value = 1.91630711295098
$ws.Range("B25").Value = 1.95259594780407
$ws.Range("C25").Value = 0.1658778482944285
$ws.Range("D25").Value = 0.08682221938824952
$ws.Range("F25").Value = 2.324136110445693
$ws.Range("G25").Value = 1.680203399650054
$ws.Range("H25").Value = 1.46510859146133
$ws.Range("J25").Value = 0.1830995110332974
$ws.Range("L25").Value = 0.3321264115533467
$ws.Range("M25").Value = 0.4325629995732214
$ws.Range("N25").Value = 1.95557339875884
